$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 87.05896960776526
$ws.Range("C2").Value = 121.5908820730844
$ws.Range("D2").Value = 139.1267513772549
$ws.Range("E2").Value = 151.0904561212856

$ws.Range("B3").Value = 108.4127932275976
$ws.Range("C3").Value = 150.7408927521585
$ws.Range("D3").Value = 170.1596471166625
$ws.Range("E3").Value = 184.449893567176

$ws.Range("B4").Value = 88.11090091047247
$ws.Range("C4").Value = 125.3638410699251
$ws.Range("D4").Value = 144.9366791515159
$ws.Range("E4").Value = 161.4687760621436

$ws.Range("B5").Value = 77.31666034521888
$ws.Range("C5").Value = 106.4939872810278
$ws.Range("D5").Value = 117.4312684384487
$ws.Range("E5").Value = 127.0904548670038

$ws.Range("B6").Value = 67.31403954872444
$ws.Range("C6").Value = 92.36376584674709
$ws.Range("D6").Value = 102.7773699613961
$ws.Range("E6").Value = 110.145913017866

$ws.Range("B7").Value = 7.417616878112186
$ws.Range("C7").Value = 10.09324119835396
$ws.Range("D7").Value = 11.2473686372922
$ws.Range("E7").Value = 11.90631588347215

$ws.Range("B8").Value = 354.1652498938367
$ws.Range("C8").Value = 494.5730662335637
$ws.Range("D8").Value = 568.1351722164003
$ws.Range("E8").Value = 607.7426005555698

$ws.Range("B9").Value = 102.6346280693362
$ws.Range("C9").Value = 140.6306409278736
$ws.Range("D9").Value = 155.9050442419406
$ws.Range("E9").Value = 165.6749734071279

$ws.Range("B10").Value = 44.57510951385223
$ws.Range("C10").Value = 59.11961234769119
$ws.Range("D10").Value = 65.45753887596777
$ws.Range("E10").Value = 67.23620524088079

$ws.Range("B11").Value = 8.121063019364485
$ws.Range("C11").Value = 10.20766562908081
$ws.Range("D11").Value = 11.17089555135803
$ws.Range("E11").Value = 12.38720678934606

$ws.Range("B12").Value = 18.46051038324407
$ws.Range("C12").Value = 24.6250343134057
$ws.Range("D12").Value = 26.04882697309632
$ws.Range("E12").Value = 25.63671716669452

$ws.Range("B13").Value = 24.99308611141655
$ws.Range("C13").Value = 32.81894819064438
$ws.Range("D13").Value = 36.64986214925756
$ws.Range("E13").Value = 38.07595649708872
